$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New QUINA lottery draw results (contest 6793-6827) appended below the
# existing history table, replacing the old last row (6793) which had
# stale numbers.
$rows = @(
    @(6793,1,24,30,62,71),
    @(6794,3,17,37,62,76),
    @(6795,1,10,15,39,69),
    @(6796,5,17,18,20,45),
    @(6797,11,21,42,50,60),
    @(6798,4,19,32,37,43),
    @(6799,2,6,11,55,76),
    @(6800,14,16,18,29,32),
    @(6801,5,10,31,75,76),
    @(6802,27,33,43,65,71),
    @(6803,6,47,50,54,58),
    @(6804,6,16,18,26,73),
    @(6805,1,3,14,19,56),
    @(6806,3,17,27,55,70),
    @(6807,7,19,21,37,65),
    @(6808,3,16,38,48,66),
    @(6809,7,8,34,61,75),
    @(6810,10,17,33,49,60),
    @(6811,5,12,21,50,57),
    @(6812,2,14,38,39,52),
    @(6813,12,36,38,42,80),
    @(6814,4,14,26,66,74),
    @(6815,1,27,50,71,72),
    @(6816,2,13,52,60,63),
    @(6817,6,13,36,53,63),
    @(6818,6,38,43,52,63),
    @(6819,1,11,32,37,62),
    @(6820,15,41,59,72,79),
    @(6821,1,2,27,62,79),
    @(6822,3,32,49,50,52),
    @(6823,1,12,29,70,74),
    @(6824,4,11,12,29,37),
    @(6825,7,9,10,34,70),
    @(6826,1,52,53,66,78),
    @(6827,1,52,53,66,78)
)

$startRow = 302
$numRows = $rows.Count
$numCols = 6
$data = New-Object 'object[,]' $numRows, $numCols
for ($i = 0; $i -lt $numRows; $i++) {
    for ($j = 0; $j -lt $numCols; $j++) {
        $data[$i, $j] = $rows[$i][$j]
    }
}

$endRow = $startRow + $numRows - 1
$ws.Range("A$startRow`:F$endRow").Value = $data

# Update the sheet dimension/view to match where Excel left the cursor
# after entering the new rows (last row fully selected).
$excel.ActiveWindow.ScrollRow = 299
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("A302:XFD302").Select()
